$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of performance test data (rows 27-30)
$ws.Range("A27").Value = "RegistrarVeterinario"
$ws.Range("B27").Value = 5000
$ws.Range("C27").Value = 7000

$ws.Range("A28").Value = "EditarVeterinario"
$ws.Range("B28").Value = 4000
$ws.Range("C28").Value = 7000

$ws.Range("A29").Value = "RegistarVisita"
$ws.Range("B29").Value = 150
$ws.Range("C29").Value = 1000
$ws.Range("E29").Value = "Profiling"

$ws.Range("A30").Value = "EditarVisita"
$ws.Range("B30").Value = 150
$ws.Range("C30").Value = 1000
$ws.Range("E30").Value = "Profiling"

# Apply the same row styling (fill colors) used by the existing data rows
$ws.Range("A22:C22").Copy() | Out-Null
$ws.Range("A27:C28").PasteSpecial(-4122) | Out-Null

$ws.Range("A22:C22").Copy() | Out-Null
$ws.Range("A29:C30").PasteSpecial(-4122) | Out-Null

$ws.Range("D19").Select()
